$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = -1.55
$ws.Range("E2").Value = -0.8

$ws.Range("C3").Value = -1.2000000000000002
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = -2.0500000000000003

[void]$ws.Range("B1:E3").Select()
